# grip_instruction.pptx — "remove getinfo explicit function"
#
# 1) The slide-master / every slide-layout footer contains a
#    datetimeFigureOut field whose cached text is "2019/6/25"; bump it to
#    "2019/6/30" (12 occurrences: 1 master + 11 layouts).
# 2) On slide 3, shape "矩形 1" has a long instructional paragraph that is
#    split into many same-styled adjacent runs; several pairs/trios of
#    neighbouring runs (identical rPr, no visible text change) get merged
#    into a single run each.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer date placeholders on the slide master and every layout.
# ---------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "2019/6/25") {
                $sh.TextFrame.TextRange.Text = "2019/6/30"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateField $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 3 — merge adjacent identically-formatted runs in the
#    instructions paragraph (3rd paragraph of the "矩形 1" textbox).
#    The visible text is unchanged; only run boundaries collapse, so we
#    rewrite small Characters() sub-ranges that each span exactly the
#    pair/trio of runs being merged.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$box = $slide3.Shapes.Item(1)
$tr = $box.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$base = $para.Start

# (offset within paragraph, length, merged text)
$mergeSpans = @(
    @(30, 5, "指示点变为"),
    @(51, 4, "开始快速"),
    @(110, 4, "其到达下"),
    @(180, 9, "个节拍，即使到达门"),
    @(189, 22, "内，也记为错误。第五道门为终点，到达即可。每"),
    @(237, 12, "秒播放声音展示节拍快慢。")
)

foreach ($span in $mergeSpans) {
    $offset = $span[0]
    $len = $span[1]
    $text = $span[2]
    $sub = $tr.Characters($base + $offset, $len)
    $sub.Text = $text
}
